$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 166349
$ws.Range("C4").Value = 157259
$ws.Range("C5").Value = 9090
$ws.Range("C8").Value = 65.2
